$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix bugs: correct full names for two users
$ws.Range("A3").Value = "Sebastian Jerezano"
$ws.Range("A7").Value = "Manuel"
